# Auto-generated script applying the cryptos.xlsx price/volume update
# (GitHub Actions data refresh commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force the cell to remain plain text (matches the source data, which
    # stores prices/percentages as strings, not numbers) while leaving the
    # cell style untouched (no explicit style index in the original file).
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextCell 'D2' '21.168.34'
Set-TextCell 'E2' '  +3.76%  '
Set-TextCell 'D3' '1.541.21'
Set-TextCell 'E3' '  +5.50%  '
Set-TextCell 'D4' '1.009'
Set-TextCell 'E4' '  -0.04%  '
Set-TextCell 'D5' '0.9626'
Set-TextCell 'E5' '  +1.32%  '
Set-TextCell 'D6' '282.51'
Set-TextCell 'E6' '  +2.80%  '
Set-TextCell 'D7' '0.3631'
Set-TextCell 'E7' '  -0.63%  '
Set-TextCell 'D8' '0.3190'
Set-TextCell 'E8' '  +4.09%  '
Set-TextCell 'D9' '40.69'
Set-TextCell 'E9' '  +2.30%  '
Set-TextCell 'D10' '1.094'
Set-TextCell 'E10' '  +5.83%  '
Set-TextCell 'D11' '0.06818'
Set-TextCell 'E11' '  +3.70%  '
Set-TextCell 'D12' '1.003'
Set-TextCell 'E12' '  +0.19%  '
Set-TextCell 'D13' '5.676'
Set-TextCell 'E13' '  +4.75%  '
Set-TextCell 'D14' '18.74'
Set-TextCell 'E14' '  +4.70%  '
Set-TextCell 'D15' '6.358'
Set-TextCell 'E15' '  +3.62%  '
Set-TextCell 'D16' '0.00001048'
Set-TextCell 'E16' '  +2.48%  '
Set-TextCell 'B17' 'Dai'
Set-TextCell 'C17' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 'D17' '0.9616'
Set-TextCell 'E17' '  -0.63%  '
Set-TextCell 'B18' 'WrappedEther'
Set-TextCell 'C18' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 'D18' '1.531.58'
Set-TextCell 'E18' '  +4.97%  '
Set-TextCell 'D19' '0.06081'
Set-TextCell 'E19' '  +4.78%  '
Set-TextCell 'D20' '72.56'
Set-TextCell 'E20' '  +4.54%  '
Set-TextCell 'D21' '5.698'
Set-TextCell 'E21' '  +4.91%  '
Set-TextCell 'D22' '15.02'
Set-TextCell 'E22' '  +3.91%  '
Set-TextCell 'E23' '  +4.56%  '
Set-TextCell 'D24' '2.311'
Set-TextCell 'E24' '  +2.89%  '
Set-TextCell 'D25' '21.183.30'
Set-TextCell 'E25' '  +3.72%  '
Set-TextCell 'D26' '148.45'
Set-TextCell 'E26' '  +4.97%  '
Set-TextCell 'D27' '2.212'
Set-TextCell 'E27' '  +6.36%  '
Set-TextCell 'D28' '17.73'
Set-TextCell 'E28' '  +3.61%  '
Set-TextCell 'D29' '1.711.67'
Set-TextCell 'E29' '  +6.17%  '
Set-TextCell 'D30' '118.05'
Set-TextCell 'E30' '  +5.20%  '
Set-TextCell 'D31' '4.018'
Set-TextCell 'E31' '  +5.23%  '
Set-TextCell 'B32' 'ImmutableX'
Set-TextCell 'C32' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D32' '0.8531'
Set-TextCell 'E32' '  +8.03%  '
Set-TextCell 'B33' 'Filecoin'
Set-TextCell 'C33' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D33' '5.201'
Set-TextCell 'E33' '  +6.41%  '
Set-TextCell 'D34' '0.08009'
Set-TextCell 'E34' '  +1.53%  '
Set-TextCell 'D35' '1.510'
Set-TextCell 'E35' '  -1.44%  '
Set-TextCell 'D36' '4.989'
Set-TextCell 'E36' '  +6.63%  '
Set-TextCell 'D37' '1.209'
Set-TextCell 'E37' '  +4.99%  '
Set-TextCell 'D38' '0.05858'
Set-TextCell 'E38' '  +2.36%  '
Set-TextCell 'D39' '0.02113'
Set-TextCell 'E39' '  +4.49%  '
Set-TextCell 'D40' '10.68'
Set-TextCell 'E40' '  +3.55%  '
Set-TextCell 'D41' '7.744'
Set-TextCell 'E41' '  +3.62%  '
Set-TextCell 'D42' '0.1922'
Set-TextCell 'D43' '0.9608'
Set-TextCell 'E43' '  +0.47%  '
Set-TextCell 'D44' '0.5463'
Set-TextCell 'E44' '  +3.96%  '
Set-TextCell 'D45' '12.51'
Set-TextCell 'E45' '  +5.30%  '
Set-TextCell 'D46' '3.584'
Set-TextCell 'E46' '  +2.71%  '
Set-TextCell 'D47' '0.5461'
Set-TextCell 'E47' '  +6.53%  '
Set-TextCell 'D48' '121.83'
Set-TextCell 'E48' '  +4.19%  '
Set-TextCell 'D49' '1.871'
Set-TextCell 'E49' '  +7.11%  '
Set-TextCell 'D50' '0.06599'
Set-TextCell 'E50' '  +2.89%  '
Set-TextCell 'D51' '69.89'
Set-TextCell 'E51' '  +5.68%  '
